$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($range, [string]$value)
    # Force a plain-text write so date-shaped strings (e.g. "2022-05-13")
    # are not auto-coerced into an Excel date serial number.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# --- Rename sheets: delimiter changed from "_" to "." ---
$wb.Worksheets.Item("dcim_device-roles").Name = "dcim.device-roles"
$wb.Worksheets.Item("dcim_device-types").Name = "dcim.device-types"
$wb.Worksheets.Item("dcim_manufacturers").Name = "dcim.manufacturers"
$wb.Worksheets.Item("dcim_regions").Name = "dcim.regions"
$wb.Worksheets.Item("dcim_sites").Name = "dcim.sites"
$wb.Worksheets.Item("users_tokens").Name = "users.tokens"
$wb.Worksheets.Item("users_users").Name = "users.users"

# --- dcim.regions sheet: refreshed ids / urls / timestamps ---
$wsRegions = $wb.Worksheets.Item("dcim.regions")

$wsRegions.Range("A2").Value = "f978e5c3-3836-4a26-bf9b-1f04b7a6e4ae"
$wsRegions.Range("B2").Value = "http://192.168.248.144:8080/api/dcim/regions/f978e5c3-3836-4a26-bf9b-1f04b7a6e4ae/"
Set-TextValue $wsRegions.Range("J2") "2022-05-13"
$wsRegions.Range("K2").Value = "2022-05-13T04:25:40.319031Z"

$wsRegions.Range("A3").Value = "271c8d66-8f2d-4bf0-85f8-62b3989afe11"
$wsRegions.Range("B3").Value = "http://192.168.248.144:8080/api/dcim/regions/271c8d66-8f2d-4bf0-85f8-62b3989afe11/"
Set-TextValue $wsRegions.Range("J3") "2022-05-13"
$wsRegions.Range("K3").Value = "2022-05-13T04:25:40.331620Z"

$wsRegions.Range("A4").Value = "21b0dcc7-4c24-489e-8f82-4d8cdcfdd545"
$wsRegions.Range("B4").Value = "http://192.168.248.144:8080/api/dcim/regions/21b0dcc7-4c24-489e-8f82-4d8cdcfdd545/"
Set-TextValue $wsRegions.Range("J4") "2022-05-13"
$wsRegions.Range("K4").Value = "2022-05-13T04:25:40.349150Z"

$wsRegions.Range("A5").Value = "43832989-d85a-459e-9ce3-741407cd76c6"
$wsRegions.Range("B5").Value = "http://192.168.248.144:8080/api/dcim/regions/43832989-d85a-459e-9ce3-741407cd76c6/"
Set-TextValue $wsRegions.Range("J5") "2022-05-13"
$wsRegions.Range("K5").Value = "2022-05-13T04:25:40.339997Z"

$wsRegions.Range("A6").Value = "20b54a29-6268-4e79-859d-a54d60fe4c95"
$wsRegions.Range("B6").Value = "http://192.168.248.144:8080/api/dcim/regions/20b54a29-6268-4e79-859d-a54d60fe4c95/"
Set-TextValue $wsRegions.Range("J6") "2022-05-13"
$wsRegions.Range("K6").Value = "2022-05-13T04:25:40.357922Z"

# --- dcim.sites sheet: refreshed ids / urls / region refs / timestamps ---
$wsSites = $wb.Worksheets.Item("dcim.sites")

$wsSites.Range("A2").Value = "fb29f887-28be-490f-b935-dd3ff3b2c81e"
$wsSites.Range("B2").Value = "http://192.168.248.144:8080/api/dcim/sites/fb29f887-28be-490f-b935-dd3ff3b2c81e/"
$wsSites.Range("F2").Value = "f978e5c3-3836-4a26-bf9b-1f04b7a6e4ae"
Set-TextValue $wsSites.Range("V2") "2022-05-13"
$wsSites.Range("W2").Value = "2022-05-13T04:25:40.420468Z"

$wsSites.Range("A3").Value = "6064d233-ed27-4436-b964-f676d47765d2"
$wsSites.Range("B3").Value = "http://192.168.248.144:8080/api/dcim/sites/6064d233-ed27-4436-b964-f676d47765d2/"
$wsSites.Range("F3").Value = "f978e5c3-3836-4a26-bf9b-1f04b7a6e4ae"
Set-TextValue $wsSites.Range("V3") "2022-05-13"
$wsSites.Range("W3").Value = "2022-05-13T04:25:40.433948Z"

$wsSites.Range("A4").Value = "c75e17a0-6b4a-4faf-9ce1-65c9967e825a"
$wsSites.Range("B4").Value = "http://192.168.248.144:8080/api/dcim/sites/c75e17a0-6b4a-4faf-9ce1-65c9967e825a/"
$wsSites.Range("F4").Value = "f978e5c3-3836-4a26-bf9b-1f04b7a6e4ae"
Set-TextValue $wsSites.Range("V4") "2022-05-13"
$wsSites.Range("W4").Value = "2022-05-13T04:25:40.444607Z"

$wsSites.Range("A5").Value = "a0fba41c-ee34-4b94-a92b-518867e4ec50"
$wsSites.Range("B5").Value = "http://192.168.248.144:8080/api/dcim/sites/a0fba41c-ee34-4b94-a92b-518867e4ec50/"
$wsSites.Range("F5").Value = "f978e5c3-3836-4a26-bf9b-1f04b7a6e4ae"
Set-TextValue $wsSites.Range("V5") "2022-05-13"
$wsSites.Range("W5").Value = "2022-05-13T04:25:40.455912Z"

Write-Host "Edits applied"
